# Updated symbol list on Mon Dec 12 03:58:42 UTC 2022 with GitHub Actions
#
# Refresh of the "Price" column (D) with newer quotes, plus rows 42/43
# (BKEXToken / CEJI) swapping rank position in the scraped source list.
#
# NOTE: every value in this sheet is stored as TEXT (numeric-looking
# strings included). Writing a numeric-looking literal straight into
# Range.Value would let Excel auto-recast it as a Number, so each cell is
# briefly switched to the "Text" number format before the write and the
# formatting is cleared again right after -- this mirrors how the source
# file already looks (General-formatted cells holding text) while still
# forcing Excel to keep the literal as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "D2"  "283.24"
Set-TextValue "D3"  "20.55"
Set-TextValue "D4"  "6.226"
Set-TextValue "D5"  "0.06176"
Set-TextValue "D6"  "3.582"
Set-TextValue "D7"  "6.546"
Set-TextValue "D8"  "1.499"
Set-TextValue "D9"  "0.8189"
Set-TextValue "D11" "0.1626"
Set-TextValue "D12" "0.08459"
Set-TextValue "D13" "0.03475"
Set-TextValue "D14" "0.03208"
Set-TextValue "D15" "0.09188"
Set-TextValue "D16" "3.705"
Set-TextValue "D17" "0.001648"
Set-TextValue "D18" "0.04724"
Set-TextValue "D19" "0.006418"
Set-TextValue "D20" "0.006172"
Set-TextValue "D23" "3.829"
Set-TextValue "D40" "0.04714"
Set-TextValue "D41" "0.007197"

# Row 42 becomes CEJI, row 43 becomes BKEXToken (positions swapped).
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004509"
Set-TextValue "E42" "41CEJICEJI"

Set-TextValue "B43" "BKEXToken"
Set-TextValue "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1098"
Set-TextValue "E43" "42BKEXTokenBKK"

Set-TextValue "D45" "0.00006884"
Set-TextValue "D48" "0.002876"
Set-TextValue "D50" "0.01243"
